# ---------------------------------------------------------------------------
# "List of cards.xlsx" edit: reword mandatory-activation-condition cards to
# use an optional "counter" mechanic instead, add a new blank "Ideas for
# cards" sheet, and add a new empty "Sheet5".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Workbook structure: insert "Ideas for cards" after "Simple 2player list"
#    (before "Master List"), and append a new empty "Sheet5" at the end.
# ---------------------------------------------------------------------------
$simple2p = $wb.Worksheets.Item("Simple 2player list")
$ideas = $wb.Worksheets.Add($null, $simple2p)
$ideas.Name = "Ideas for cards"

$lastIdx = $wb.Worksheets.Count
$lastWs = $wb.Worksheets.Item($lastIdx)
$sheet5 = $wb.Worksheets.Add($null, $lastWs)
$sheet5.Name = "Sheet5"

# ---------------------------------------------------------------------------
# 2) "Ideas for cards": same header row as "Simple 2player list" row 1.
# ---------------------------------------------------------------------------
$ideas.Range("A1:N1").WrapText = $true
$ideas.Rows.Item(1).RowHeight = 38
$ideas.Range("A1").Value = "Name"
$ideas.Range("B1").Value = "Archetype"
$ideas.Range("C1").Value = "Class"
$ideas.Range("D1").Value = "Card Type"
$ideas.Range("E1").Value = "Attribute"
$ideas.Range("G1").Value = "Effect"
$ideas.Range("H1").Value = "Points"
$ideas.Range("K1").Value = "notes"
$ideas.Range("L1").Value = "purpose"
$ideas.Range("M1").Value = "intended function"
$ideas.Range("N1").Value = """effective"" points"
$ideas.Range("B2").Select()

# ---------------------------------------------------------------------------
# 3) "Simple 2player list": reword the three mandatory-condition cards to
#    the new optional "counter" wording, drop the explicit
#    "(This card can only be played ...)" flavour-text rows, and update the
#    points shown for the now-optional effects to 0 (the bonus is no longer
#    guaranteed).
# ---------------------------------------------------------------------------
$s2p = $wb.Worksheets.Item("Simple 2player list")

# Row 8 - Gambler: "flip a coin" effect now leaves an optional counter.
$s2p.Range("G8").Value = "When played: Flip a coin. If heads, put a counter on this card. This card gains +500 points if it has a counter on it."
$s2p.Rows.Item(8).RowHeight = 228

# Row 12 - First in Line: drop the mandatory "(you control no cards)" line,
# replace with optional counter wording; the guaranteed 1500 bonus is removed.
$s2p.Range("G12").Value = "When played: If you control no cards, place a counter on this card. If this card has a counter on it, it has +1500 points."
$s2p.Range("H12").Value = 0
$s2p.Rows.Item(12).RowHeight = 247

# Row 11 - Second in Line: drop the mandatory "(all players control...)" line,
# replace with optional counter wording; the guaranteed 700 bonus is removed.
$s2p.Range("G11").Value = "When played: (If each player controls at least 1 card) place a counter on this card. If this card has a counter on it, it thas +700 points."
$s2p.Range("H11").Value = 0
$s2p.Rows.Item(11).RowHeight = 266

# Final view state for this sheet.
$s2p.Activate()
$s2p.Range("H12").Select()
